$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the manual threshold value (B21): 225 -> 175
$ws.Range("B21").Value = 175

# Re-enter the formulas over the full ranges at once so Excel stores them
# as shared formulas (matching how Excel normally fills a block selection).
$ws.Range("A11:D14").FormulaR1C1 = "=(R[-6]C[7]+R[-6]C[12]+R[-6]C[17])/3"
$ws.Range("G11:J14").FormulaR1C1 = "=ROUND(RC[-6],0)"

# Update the sheet view (scroll position / selection) to match the saved state
$ws.Application.ActiveWindow.ScrollRow = 5
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F22").Select()
